$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of this week's data block (old rows 238-247
# shift down to 240-249), matching how this "weekly" price sheet logs the
# newest observations above older ones.
$ws.Rows("238:239").Insert()

# New row 238: Mango, "Primera" grade, from Peru.
$ws.Cells.Item(238, 1).Value = 3
$ws.Cells.Item(238, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(238, 3).Value = "Coquimbo"
$ws.Cells.Item(238, 4).Value = 44509
$ws.Cells.Item(238, 5).Value = 5
$ws.Cells.Item(238, 6).Value = "Fruta"
$ws.Cells.Item(238, 7).Value = 100108
$ws.Cells.Item(238, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(238, 9).Value = 100108002
$ws.Cells.Item(238, 10).Value = "Mango"
$ws.Cells.Item(238, 11).Value = "Sin especificar"
$ws.Cells.Item(238, 12).Value = "Primera"
$ws.Cells.Item(238, 13).Value = 228
$ws.Cells.Item(238, 14).Value = 7000
$ws.Cells.Item(238, 15).Value = 7000
$ws.Cells.Item(238, 16).Value = 7000
$ws.Cells.Item(238, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(238, 18).Value = "Perú"
$ws.Cells.Item(238, 19).Value = 1750
$ws.Cells.Item(238, 20).Value = 4

# New row 239: same date/market, "Segunda" grade, from Peru.
$ws.Cells.Item(239, 1).Value = 3
$ws.Cells.Item(239, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(239, 3).Value = "Coquimbo"
$ws.Cells.Item(239, 4).Value = 44509
$ws.Cells.Item(239, 5).Value = 5
$ws.Cells.Item(239, 6).Value = "Fruta"
$ws.Cells.Item(239, 7).Value = 100108
$ws.Cells.Item(239, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(239, 9).Value = 100108002
$ws.Cells.Item(239, 10).Value = "Mango"
$ws.Cells.Item(239, 11).Value = "Sin especificar"
$ws.Cells.Item(239, 12).Value = "Segunda"
$ws.Cells.Item(239, 13).Value = 228
$ws.Cells.Item(239, 14).Value = 7000
$ws.Cells.Item(239, 15).Value = 7000
$ws.Cells.Item(239, 16).Value = 7000
$ws.Cells.Item(239, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(239, 18).Value = "Perú"
$ws.Cells.Item(239, 19).Value = 1750
$ws.Cells.Item(239, 20).Value = 4
